# Append one new data-log row to each of the 4 worksheets.
#
# Sheet "ROW50-FE-LIFTER"  : new row 38
# Sheet "ROW50-MID-LIFTER" : new row 40
# Sheet "ROW11-FE-LIFTER"  : new row 38
# Sheet "ROW11-MID-LIFTER" : new row 38
#
# Columns: A time (date serial), B total-len raw, C ID raw, D actual-len
# raw, E checksum raw, F total-len dec, G ID dec, H actual-len dec,
# I checksum dec.

$wb = $excel.ActiveWorkbook

# The decoded-ID value (column G) is the same huge integer
# (568631262647113771663628) everywhere in this workbook. It doesn't fit
# exactly in a double, so most sheets store Excel's rounded double
# (5.68631262647114e+23); build that by multiplying (this engine's
# PowerShell parser does not accept scientific-notation literals, nor
# integer literals too big for a native int type).
$idDecNumber = 568631262647114 * 1000000000

# ---------------- Sheet 1: ROW50-FE-LIFTER, new row 38 ----------------
$ws = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$r = 38
$ws.Cells.Item($r, 1).Value = 45743.16792305555
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x6e"
$ws.Cells.Item($r, 5).Value = "0xe"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = $idDecNumber
$ws.Cells.Item($r, 8).Value = 366
$ws.Cells.Item($r, 9).Value = 14

# ------------- Sheet 2: ROW50-MID-LIFTER, new row 40 -------------
# Note: B has a trailing space, and G keeps the exact integer as TEXT
# (matches how this sheet already stores the same ID_DEC value elsewhere).
$ws = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$r = 40
$ws.Cells.Item($r, 1).Value = 45743.1370949074
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x90 "
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x76"
$ws.Cells.Item($r, 5).Value = "0x19"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).NumberFormat = "@"
$ws.Cells.Item($r, 7).Value = "568631262647113771663628"
$ws.Cells.Item($r, 7).Style = "Normal"
$ws.Cells.Item($r, 8).Value = 374
$ws.Cells.Item($r, 9).Value = 25

# ---------------- Sheet 3: ROW11-FE-LIFTER, new row 38 ----------------
$ws = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$r = 38
$ws.Cells.Item($r, 1).Value = 45743.18440581019
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x6e"
$ws.Cells.Item($r, 5).Value = "0x14"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = $idDecNumber
$ws.Cells.Item($r, 8).Value = 366
$ws.Cells.Item($r, 9).Value = 20

# ---------------- Sheet 4: ROW11-MID-LIFTER, new row 38 ----------------
$ws = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$r = 38
$ws.Cells.Item($r, 1).Value = 45743.33174475694
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x76"
$ws.Cells.Item($r, 5).Value = "0x19"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = $idDecNumber
$ws.Cells.Item($r, 8).Value = 374
$ws.Cells.Item($r, 9).Value = 25
